# Updated cryptos list (GitHub Actions scrape refresh): prices + 1h volume
# change %, plus a TheGraph/ONDO row-order swap (rows 50-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells that look numeric get NumberFormat "@" first so Excel keeps
# them as literal text (preserving things like trailing zeros / "1.00")
# instead of auto-coercing the Range.Value assignment to a number.

$ws.Range("D2").Value = "71.083.44"
$ws.Range("E2").Value = "  +2.89%  "
$ws.Range("D3").Value = "3.805.85"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "711.48"
$ws.Range("E5").Value = "  +13.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.23"
$ws.Range("E6").Value = "  +4.74%  "
$ws.Range("D7").Value = "3.805.30"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +3.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.47"
$ws.Range("E11").Value = "  +9.60%  "
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  +9.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.34"
$ws.Range("E14").Value = "  +4.26%  "
$ws.Range("D15").Value = "4.446.69"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "3.805.41"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "71.126.51"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.91"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.23"
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.20"
$ws.Range("E21").Value = "  +17.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.68"
$ws.Range("E22").Value = "  +3.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.717"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000148"
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.97"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.43"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.55"
$ws.Range("E27").Value = "  +4.23%  "
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("D29").Value = "3.956.73"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.17"
$ws.Range("E30").Value = "  +18.89%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  +6.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.29"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.179"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  +4.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "3.756.21"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.56"
$ws.Range("E40").Value = "  +9.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.00"
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("E42").Value = "  +11.74%  "
$ws.Range("E43").Value = "  +26.66%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "162.54"
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.45"
$ws.Range("E48").Value = "  +5.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.91"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.301"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.39"
$ws.Range("E51").Value = "  -1.69%  "
